# daily auto push: 2025-10-01 02:16 UTC
# Append the two new daily log rows to the bottom of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: 2025/10/01, 水, 9, 161
# Force column A to be treated as text so the date-like string "2025/10/01"
# is stored literally instead of being auto-converted into a date serial
# number by Excel's cell-value type inference.
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "2025/10/01"
$ws.Range("B44").Value = "水"
$ws.Range("C44").Value = 9
$ws.Range("D44").Value = 161
$ws.Range("A44").ClearFormats()

# Row 45: 2025/10/01, 水, 11, 160
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "2025/10/01"
$ws.Range("B45").Value = "水"
$ws.Range("C45").Value = 11
$ws.Range("D45").Value = 160
$ws.Range("A45").ClearFormats()
